# Update countries & provincias Spain
# Refresh the "Pais" sheet with the 27-Mar-2020 13:12 data snapshot:
#  - updates the "last updated" timestamp in A1
#  - refreshes case totals for several countries whose numbers changed
#  - a few countries swapped ranking order (by total cases), so both the
#    country name and its row of figures move together

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1
$ws.Cells.Item(1, 1).Value = "Datos actualizados a 27 de Marzo de 2020 a las 13:12"

# Row 15
$ws.Cells.Item(15, 1).Value = "Austria"
$ws.Cells.Item(15, 2).Value = 7317
$ws.Cells.Item(15, 3).Value = 408
$ws.Cells.Item(15, 4).Value = 225
$ws.Cells.Item(15, 5).Value = 7034
$ws.Cells.Item(15, 6).Value = 128
$ws.Cells.Item(15, 7).Value = 9
$ws.Cells.Item(15, 8).Value = 58

# Row 16
$ws.Cells.Item(16, 1).Value = "Belgica"
$ws.Cells.Item(16, 2).Value = 7284
$ws.Cells.Item(16, 3).Value = 1049
$ws.Cells.Item(16, 4).Value = 858
$ws.Cells.Item(16, 5).Value = 6137
$ws.Cells.Item(16, 6).Value = 690
$ws.Cells.Item(16, 7).Value = 69
$ws.Cells.Item(16, 8).Value = 289

# Row 24
$ws.Cells.Item(24, 2).Value = 2893
$ws.Cells.Item(24, 3).Value = 53
$ws.Cells.Item(24, 5).Value = 2800
$ws.Cells.Item(24, 6).Value = 214

# Row 45
$ws.Cells.Item(45, 2).Value = 799
$ws.Cells.Item(45, 3).Value = 72
$ws.Cells.Item(45, 4).Value = 73
$ws.Cells.Item(45, 5).Value = 706

# Row 53
$ws.Cells.Item(53, 6).Value = 7

# Row 60
$ws.Cells.Item(60, 4).Value = 227
$ws.Cells.Item(60, 5).Value = 235

# Row 62
$ws.Cells.Item(62, 5).Value = 434
$ws.Cells.Item(62, 7).Value = 1
$ws.Cells.Item(62, 8).Value = 8

# Row 68
$ws.Cells.Item(68, 4).Value = 28
$ws.Cells.Item(68, 5).Value = 300

# Row 71
$ws.Cells.Item(71, 4).Value = 9
$ws.Cells.Item(71, 5).Value = 264

# Row 73
$ws.Cells.Item(73, 6).Value = 1

# Row 79
$ws.Cells.Item(79, 4).Value = 5
$ws.Cells.Item(79, 5).Value = 216

# Row 82
$ws.Cells.Item(82, 1).Value = "Republica de Macedonia"
$ws.Cells.Item(82, 2).Value = 219
$ws.Cells.Item(82, 3).Value = 18
$ws.Cells.Item(82, 4).Value = 3
$ws.Cells.Item(82, 5).Value = 213
$ws.Cells.Item(82, 6).Value = 1
$ws.Cells.Item(82, 8).Value = 3

# Row 83
$ws.Cells.Item(83, 1).Value = "Jordania"
$ws.Cells.Item(83, 2).Value = 212
$ws.Cells.Item(83, 4).Value = 2
$ws.Cells.Item(83, 5).Value = 210
$ws.Cells.Item(83, 6).Value = 0
$ws.Cells.Item(83, 8).Value = 0

# Row 84
$ws.Cells.Item(84, 1).Value = "San Marino"
$ws.Cells.Item(84, 2).Value = 208
$ws.Cells.Item(84, 4).Value = 4
$ws.Cells.Item(84, 5).Value = 183
$ws.Cells.Item(84, 6).Value = 12
$ws.Cells.Item(84, 8).Value = 21

# Row 99
$ws.Cells.Item(99, 4).Value = 31
$ws.Cells.Item(99, 5).Value = 75

# Row 113
$ws.Cells.Item(113, 2).Value = 68
$ws.Cells.Item(113, 3).Value = 1
$ws.Cells.Item(113, 5).Value = 67

# Row 149
$ws.Cells.Item(149, 1).Value = "Republica de Yibuti"
$ws.Cells.Item(149, 3).Value = 1

# Row 150
$ws.Cells.Item(150, 1).Value = "Guinea Ecuatorial"
$ws.Cells.Item(150, 2).Value = 12
$ws.Cells.Item(150, 5).Value = 12

# Row 151
$ws.Cells.Item(151, 1).Value = "San Martin (Parte Francesa)"

# Row 152
$ws.Cells.Item(152, 1).Value = "Dominica"

# Row 153
$ws.Cells.Item(153, 1).Value = "Mongolia"
